$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1477845621989524
$ws.Range("C2").Value = 0.1738072731056164
$ws.Range("B3").Value = 0.1023952415260214
$ws.Range("C3").Value = 0.1667320816855214
$ws.Range("B4").Value = 0.08231148181531057
$ws.Range("C4").Value = 0.1035517964098561
$ws.Range("B5").Value = 0.08152802259139066
$ws.Range("C5").Value = 0.08644468207250083
$ws.Range("B6").Value = 0.0622064487966581
$ws.Range("C6").Value = 0.1263389879741443
